# Fix up raw data strings in the category/value key sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "planktonic larval duration (PLD) exposure"
$ws.Range("D5").Value = " <1; 1-2; 2-5; 5-10; 10-20; 20-50; 50-100; 100-1000; 1000-10,000; >10,000"
$ws.Range("C6").Value = "1, 2-10, 11-25, 26-50, 51-100, >100"
$ws.Range("D6").Value = "1; 2-10; 11-25; 26-50; 51-100; >100"
$ws.Range("D37").Value = " >1000 mm; 50 mm-999 mm; 0.5mm-49 mm; <0.4 mm"
$ws.Range("C44").Value = "high; medium; low; none"
$ws.Range("D44").Value = "high; medium; low; none"

$ws.Range("D6").Select()
